# Remove the Emilia Romagna Grand Prix rows from the F1 schedule sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("F1")

# Rows 40-44 hold the 5 sessions (FP1, FP2, FP3, Qualifying, Race) for the
# Emilia Romagna Grand Prix. Deleting the entire rows shifts everything
# below up by 5, matching the target workbook state.
$rows = $ws.Range("A40:I44").EntireRow
$rows.Delete()

# Clear any stale multi-cell selection left over from editing (Excel resets
# the sheet view selection to a single default cell after such edits).
$ws.Range("A1").Select()
